$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '54.222.19'
$ws.Range('E2').Value = '  -7.41%  '

$ws.Range('D3').Value = '2.412.20'
$ws.Range('E3').Value = '  -11.41%  '

$ws.Range('E4').Value = '  +0.28%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '464.81'
$ws.Range('E5').Value = '  -7.68%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.30'
$ws.Range('E6').Value = '  -6.90%  '

$ws.Range('E7').Value = '  +0.11%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.490'
$ws.Range('E8').Value = '  -7.35%  '

$ws.Range('D9').Value = '2.429.14'
$ws.Range('E9').Value = '  -11.09%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0948'
$ws.Range('E10').Value = '  -9.50%  '

$ws.Range('E11').Value = '  -12.35%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.313'
$ws.Range('E12').Value = '  -9.83%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.121'
$ws.Range('E13').Value = '  -4.08%  '

$ws.Range('D14').Value = '2.866.60'
$ws.Range('E14').Value = '  -10.52%  '

$ws.Range('D15').Value = '54.272.58'
$ws.Range('E15').Value = '  -7.46%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000132'
$ws.Range('E16').Value = '  -2.23%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '19.67'
$ws.Range('E17').Value = '  -9.38%  '

$ws.Range('D18').Value = '2.447.29'
$ws.Range('E18').Value = '  -9.90%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.18'
$ws.Range('E19').Value = '  -12.30%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '309.74'
$ws.Range('E20').Value = '  -9.91%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.43'
$ws.Range('E21').Value = '  -13.93%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.993'
$ws.Range('E22').Value = '  -0.45%  '

$ws.Range('E23').Value = '  +0.10%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.36'
$ws.Range('E24').Value = '  -14.36%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '56.09'
$ws.Range('E25').Value = '  -10.76%  '

$ws.Range('E26').Value = '  +1.15%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.384'
$ws.Range('E27').Value = '  -10.13%  '

$ws.Range('E28').Value = '  -10.09%  '

$ws.Range('D29').Value = '2.533.75'
$ws.Range('E29').Value = '  -10.86%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.11'
$ws.Range('E30').Value = '  -5.24%  '

$ws.Range('E31').Value = '  +0.19%  '

$ws.Range('D32').Value = '0.0₃0711'
$ws.Range('E32').Value = '  -14.49%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '146.75'
$ws.Range('E33').Value = '  -2.79%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '17.70'
$ws.Range('E34').Value = '  -7.59%  '

$ws.Range('E35').Value = '  -10.71%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.99'
$ws.Range('E36').Value = '  -8.10%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.53'
$ws.Range('E37').Value = '  -15.96%  '

$ws.Range('E38').Value = '  -7.29%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.797'
$ws.Range('E39').Value = '  -16.10%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.995'
$ws.Range('E40').Value = '  -0.07%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '32.82'
$ws.Range('E41').Value = '  -8.57%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.594'
$ws.Range('E42').Value = '  -1.12%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0521'
$ws.Range('E43').Value = '  -6.78%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.23'
$ws.Range('E44').Value = '  -8.93%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.23'
$ws.Range('E46').Value = '  -11.84%  '

$ws.Range('D47').Value = '1.935.32'
$ws.Range('E47').Value = '  -11.44%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0877'
$ws.Range('E48').Value = '  -0.84%  '

$ws.Range('E49').Value = '  -4.80%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '16.45'
$ws.Range('E50').Value = '  -13.49%  '

$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.16'
$ws.Range('E51').Value = '  -13.37%  '
